# soldier_card.xlsx: extend the cs-template table with a new "desc" field
# (id/name/logic_res/tex/heal -> + desc), mirroring the existing
# name/type/comment header rows, so the C# array template ("role card")
# generator picks it up as a new column.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 1: field name, Row 2: field type, Row 3: Chinese comment/description.
$ws.Cells.Item(1, 6).Value = "desc"
$ws.Cells.Item(2, 6).Value = "string"
$ws.Cells.Item(3, 6).Value = "描述"

# Leave row 4 (the SOLDIER_1000 data row) untouched for the new column -
# no value has been filled in for it yet.

# Move/leave the selection on the newly added cell, as in the authored edit.
$ws.Range("F4").Select()
